$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Name" header to "CustomerName"
$ws.Range("A1").Value = "CustomerName"

# Remove the "City" column (D) entirely
$ws.Range("D1:D4").Delete()
